$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("D9").Value = "[공고] 일반 MBA/DBA 프로그램 총괄 담당 교수 채용"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-general-mba-dba-chair/#utm_source=rss&utm_medium=rss&utm_campaign=notice-general-mba-dba-chair"

# Row 16
$ws.Range("D16").Value = "Grad-CAM++ 내용정리 [XAI-3]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/144"

# Row 32
$ws.Range("D32").Value = "Bigquery Procedure 소개"
$ws.Range("E32").Value = "https://dodonam.tistory.com/313"

# Row 36
$ws.Range("D36").Value = "Introduction to Image Super-Resolution"

# Row 42
$ws.Range("D42").Value = "MFC - typeid를 이용한 문자열 형변환 소스코드"
$ws.Range("E42").Value = "https://kjk92.tistory.com/67"

# Row 51
$ws.Range("D51").Value = "이진 분류기 성능 평가방법 AUC(area under the ROC curve)의 이해"
$ws.Range("E51").Value = "https://bskyvision.com/1165"
